$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H33").Value = 17242144
$ws_ALC.Range("I33").Value = 27778344
$ws_ALC.Range("K33").Value = 27778344
$ws_ALC.Range("M33").Value = -27778115
$ws_ALC.Range("H40").Value = 3188.276
$ws_ALC.Range("I40").Value = 1193.3636
$ws_ALC.Range("J40").Value = 4407.3887
$ws_ALC.Range("K40").Value = 1193.3636
$ws_ALC.Range("L40").Value = 4407.3887
$ws_ALC.Range("M40").Value = -1018.3636
$ws_ALC.Range("N40").Value = -4757.3887
$ws_ALC.Range("H64").Value = 3105.587
$ws_ALC.Range("I64").Value = 2923.6155
$ws_ALC.Range("K64").Value = 2923.6155
$ws_ALC.Range("M64").Value = -2675.6155
$ws_ALC.Range("H67").Value = 3105.587
$ws_ALC.Range("I67").Value = 2923.6155
$ws_ALC.Range("K67").Value = 2923.6155
$ws_ALC.Range("M67").Value = -2065.6155
$ws_ALC.Range("H74").Value = 5165.3125
$ws_ALC.Range("I74").Value = 4666.5
$ws_ALC.Range("J74").Value = 5996.6665
$ws_ALC.Range("K74").Value = 4666.5
$ws_ALC.Range("L74").Value = 5996.6665
$ws_ALC.Range("M74").Value = -3730.5
$ws_ALC.Range("N74").Value = -7868.6665
$ws_ALC.Range("H77").Value = 5165.3125
$ws_ALC.Range("I77").Value = 4666.5
$ws_ALC.Range("J77").Value = 5996.6665
$ws_ALC.Range("K77").Value = 23332.5
$ws_ALC.Range("L77").Value = 29983.3325
$ws_ALC.Range("M77").Value = -18652.5
$ws_ALC.Range("N77").Value = -39343.3325
$ws_ALC.Range("H92").Value = 1962.0968
$ws_ALC.Range("I92").Value = 2168.5186
$ws_ALC.Range("J92").Value = 568.75
$ws_ALC.Range("K92").Value = 2168.5186
$ws_ALC.Range("L92").Value = 568.75
$ws_ALC.Range("M92").Value = -920.5185999999999
$ws_ALC.Range("N92").Value = -3064.75
$ws_ALC.Range("H98").Value = 906.73914
$ws_ALC.Range("I98").Value = 614.44446
$ws_ALC.Range("J98").Value = 1959
$ws_ALC.Range("K98").Value = 614.44446
$ws_ALC.Range("L98").Value = 1959
$ws_ALC.Range("M98").Value = 883.55554
$ws_ALC.Range("N98").Value = -4955
$ws_ALC.Range("H113").Value = 3183.4827
$ws_ALC.Range("I113").Value = 2179.3076
$ws_ALC.Range("J113").Value = 3999.375
$ws_ALC.Range("K113").Value = 2179.3076
$ws_ALC.Range("L113").Value = 3999.375
$ws_ALC.Range("M113").Value = 1074.6924
$ws_ALC.Range("N113").Value = -10507.375
$ws_ALC.Range("H122").Value = 906.73914
$ws_ALC.Range("I122").Value = 614.44446
$ws_ALC.Range("J122").Value = 1959
$ws_ALC.Range("K122").Value = 1843.33338
$ws_ALC.Range("L122").Value = 5877
$ws_ALC.Range("M122").Value = 606.66662
$ws_ALC.Range("N122").Value = -10777
$ws_ALC.Range("H132").Value = 3568.5417
$ws_ALC.Range("I132").Value = 1815.9459
$ws_ALC.Range("J132").Value = 9463.637000000001
$ws_ALC.Range("K132").Value = 5447.8377
$ws_ALC.Range("L132").Value = 28390.911
$ws_ALC.Range("M132").Value = -2917.8377
$ws_ALC.Range("N132").Value = -33450.911
$ws_ALC.Range("H137").Value = 20641.4
$ws_ALC.Range("I137").Value = 24053.912
$ws_ALC.Range("J137").Value = 3199.6667
$ws_ALC.Range("K137").Value = 72161.736
$ws_ALC.Range("L137").Value = 9599.000100000001
$ws_ALC.Range("M137").Value = -69611.736
$ws_ALC.Range("N137").Value = -14699.0001
$ws_ALC.Range("H141").Value = 1505.7142
$ws_ALC.Range("I141").Value = 1505.7142
$ws_ALC.Range("K141").Value = 4517.142599999999
$ws_ALC.Range("M141").Value = 662.8574000000008
$ws_ARM.Range("H23").Value = 9766.666999999999
$ws_ARM.Range("J23").Value = 9766.666999999999
$ws_ARM.Range("L23").Value = 9766.666999999999
$ws_ARM.Range("N23").Value = -10284.667
$ws_ARM.Range("H32").Value = 6940.636
$ws_ARM.Range("I32").Value = 5350.984
$ws_ARM.Range("J32").Value = 13511.2
$ws_ARM.Range("K32").Value = 5350.984
$ws_ARM.Range("L32").Value = 13511.2
$ws_ARM.Range("M32").Value = -5063.984
$ws_ARM.Range("N32").Value = -14085.2
$ws_ARM.Range("H63").Value = 11832.333
$ws_ARM.Range("I63").Value = 13598.8
$ws_ARM.Range("K63").Value = 13598.8
$ws_ARM.Range("M63").Value = -12912.8
$ws_ARM.Range("H66").Value = 11832.333
$ws_ARM.Range("I66").Value = 13598.8
$ws_ARM.Range("K66").Value = 67994
$ws_ARM.Range("M66").Value = -64562
$ws_BSM.Range("H86").Value = 7140.421
$ws_BSM.Range("I86").Value = 4935.5
$ws_BSM.Range("J86").Value = 18900
$ws_BSM.Range("K86").Value = 4935.5
$ws_BSM.Range("L86").Value = 18900
$ws_BSM.Range("M86").Value = -3812.5
$ws_BSM.Range("N86").Value = -21146
$ws_BSM.Range("H89").Value = 7140.421
$ws_BSM.Range("I89").Value = 4935.5
$ws_BSM.Range("J89").Value = 18900
$ws_BSM.Range("K89").Value = 24677.5
$ws_BSM.Range("L89").Value = 94500
$ws_BSM.Range("M89").Value = -19061.5
$ws_BSM.Range("N89").Value = -105732
$ws_CUL.Range("H97").Value = 658.5454999999999
$ws_CUL.Range("I97").Value = 100.8
$ws_CUL.Range("J97").Value = 1123.3334
$ws_CUL.Range("K97").Value = 302.4
$ws_CUL.Range("L97").Value = 3370.0002
$ws_CUL.Range("M97").Value = 193.6
$ws_CUL.Range("N97").Value = -4362.0002
$ws_CUL.Range("H134").Value = 2436.4092
$ws_CUL.Range("J134").Value = 3554.6
$ws_CUL.Range("L134").Value = 10663.8
$ws_CUL.Range("N134").Value = -20803.8
$ws_GSM.Range("H46").Value = 20000
$ws_GSM.Range("J46").Value = 20000
$ws_GSM.Range("L46").Value = 20000
$ws_GSM.Range("N46").Value = -20312
$ws_GSM.Range("H57").Value = 0
$ws_GSM.Range("J57").Value = 0
$ws_GSM.Range("L57").Value = 0
$ws_GSM.Range("N57").Value = ""
$ws_GSM.Range("H113").Value = 2388.8
$ws_GSM.Range("I113").Value = 2000
$ws_GSM.Range("J113").Value = 2486
$ws_GSM.Range("K113").Value = 2000
$ws_GSM.Range("L113").Value = 2486
$ws_GSM.Range("M113").Value = 170
$ws_GSM.Range("N113").Value = -6826
$ws_LTW.Range("H18").Value = 14531.692
$ws_LTW.Range("I18").Value = 4445.4546
$ws_LTW.Range("J18").Value = 70006
$ws_LTW.Range("K18").Value = 4445.4546
$ws_LTW.Range("L18").Value = 70006
$ws_LTW.Range("M18").Value = -4273.4546
$ws_LTW.Range("N18").Value = -70350
$ws_LTW.Range("H20").Value = 164280.72
$ws_LTW.Range("I20").Value = 252250
$ws_LTW.Range("J20").Value = 46988.332
$ws_LTW.Range("K20").Value = 252250
$ws_LTW.Range("L20").Value = 46988.332
$ws_LTW.Range("M20").Value = -252024
$ws_LTW.Range("N20").Value = -47440.332
$ws_WVR.Range("H136").Value = 31010092
$ws_WVR.Range("I136").Value = 37038932
$ws_WVR.Range("J136").Value = 20836428
$ws_WVR.Range("K136").Value = 111116796
$ws_WVR.Range("L136").Value = 62509284
$ws_WVR.Range("M136").Value = -111114246
$ws_WVR.Range("N136").Value = -62514384
